$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the A1 conversion note text ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 13.74 = 55702.96 pesos`n✅ 55702.96 pesos = 13.73 = 971.11 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- Sheet "tasas": update rate values ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 72.79900000000001
$ws2.Range("O10").Value = 4055.12
$ws2.Range("N12").Value = 4058
$ws2.Range("O12").Value = 70.746
